$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.065.55"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.891.24"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "0.7435"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "243.13"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "0.9982"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "0.3174"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "0.07254"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "25.02"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").Value = "0.08366"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.968.26"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7603"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "5.426"
$ws.Range("D15").Value = "92.81"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "6.163"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "30.071.11"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "249.81"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "0.000007872"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.138.44"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.9971"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "8.017"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "0.9973"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "0.1583"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "9.311"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "164.71"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "18.77"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "2.055"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("D31").Value = "4.607"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").Value = "1.538"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "4.227"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("D34").Value = "0.05384"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "1.256"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "0.7625"
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("D37").Value = "0.9993"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "2.724"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "0.01974"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "2.768"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("D42").Value = "1.102.43"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "73.13"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "6.071"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "0.8736"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "0.9994"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").Value = "7.642"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "9.593"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "2.030.41"
$ws.Range("E51").Value = "  -0.68%  "
